$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data grid (A:F, rows 1-18) after the edit.
$data = @(
    @("year", "NTS", "HKL", "NTN", "KL", "HK"),
    @(2002, 0, 0, 20, 0, 0),
    @(2003, 0, 0, 1, 0, 0),
    @(2004, 0, 0, 0, 0, 0),
    @(2005, 0, 0, 0, 0, 0),
    @(2006, 0, 0, 0, 0, 0),
    @(2007, 0, 0, 0, 0, 0),
    @(2008, 0, 0, 0, 0, 0),
    @(2009, 0, 0, 0, 0, 0),
    @(2010, 0, 4, 0, 0, 4),
    @(2011, 0, 0, 0, 0, 0),
    @(2012, 0, 0, 0, 0, 0),
    @(2013, 0, 0, 0, 0, 0),
    @(2014, 0, 2, 1, 0, 2),
    @(2015, 0, 0, 3, 0, 0),
    @(2016, 0, 4, 0, 1, 3),
    @(2017, 0, 1, 0, 1, 0),
    @(2018, 10, 18, 0, 18, 0)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Match the new selection recorded in the saved workbook.
[void]$ws.Range("E1:F18").Select()
